$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New headers for columns I and J - copy formatting (bold, border, centered) from H1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data for column I (I0) and column J (IF), rows 2-14
$i0 = @(1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 1, 5)
$if = @(2, 4, 4, 5, 4, 5, 5, 6, 6, 6, 5, 5, 5)

for ($r = 0; $r -lt 13; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
